$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix product-name capitalization
$ws.Range("B17").Value = "Botin Exotic "
$ws.Range("B19").Value = "Botin Kim "
$ws.Range("B20").Value = "Botin Santorini"

# Fix image URLs to match corrected folder/file casing
$ws.Range("O13").Value = "[""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Avejita/Botin+Avejita+1+Negro.jpg""]"
$ws.Range("O17").Value = "[""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Exotic/Botin+Exotic+1+Negro.jpg""]"
$ws.Range("O19").Value = "[""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Kim/Botin+Kim+1+Blanco.jpg"", ""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Kim/Botin+Kim+2+Negro.PNG""]"
$ws.Range("O20").Value = "[""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Santorini/Botin+Santorini+1+Beige.jpg"", ""https://recursosmolova.s3.amazonaws.com/Products+Images/Majo'R/Botin+Santorini/Botin+Santorini+2+Negro.jpg""]"
